$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- "First page" footer (footer1.xml) : Pearson logo, docPr id="3" ---
# name="image1.png" -> name="image2.png"
$ftFirst = $sec.Footers.Item(2)   # wdHeaderFooterFirstPage = 2
$shpFtFirst = $ftFirst.Range.InlineShapes.Item(1)
$rngFtFirst = $shpFtFirst.Range
$xmlFtFirst = $rngFtFirst.WordOpenXML
$xmlFtFirst = $xmlFtFirst -replace 'name="image1\.png"', 'name="image2.png"'
$rngFtFirst.WordOpenXML = $xmlFtFirst

# --- "Default/primary" footer (footer2.xml) : Pearson logo, docPr id="2" ---
# name="image1.png" -> name="image2.png"
$ftDefault = $sec.Footers.Item(1)  # wdHeaderFooterPrimary = 1
$shpFtDefault = $ftDefault.Range.InlineShapes.Item(1)
$rngFtDefault = $shpFtDefault.Range
$xmlFtDefault = $rngFtDefault.WordOpenXML
$xmlFtDefault = $xmlFtDefault -replace 'name="image1\.png"', 'name="image2.png"'
$rngFtDefault.WordOpenXML = $xmlFtDefault

# --- "First page" header (header1.xml) : BTec logo, docPr id="1" ---
# name="image2.jpg" -> name="image1.jpg"
$hdFirst = $sec.Headers.Item(2)   # wdHeaderFooterFirstPage = 2
$shpHdFirst = $hdFirst.Range.InlineShapes.Item(1)
$rngHdFirst = $shpHdFirst.Range
$xmlHdFirst = $rngHdFirst.WordOpenXML
$xmlHdFirst = $xmlHdFirst -replace 'name="image2\.jpg"', 'name="image1.jpg"'
$rngHdFirst.WordOpenXML = $xmlHdFirst
